# feat: add 2022-Q3 data
#
# Before:  Sheet1 "总计" (totals) + Sheet2 "2022-Q1" (fund holdings for 2022-Q1)
# After:   Sheet1 "总计" (totals, now with a 2022-Q3 row too)
#          Sheet2 "2022-Q3" (NEW fund holdings data for 2022-Q3)
#          Sheet3 "2022-Q1" (the old fund holdings data, relocated)

$wb = $excel.ActiveWorkbook
$totals = $wb.Worksheets.Item(1)
$q1Sheet = $wb.Worksheets.Item(2)

# A scratch cell used to force numeric-looking strings (e.g. "005535",
# "10.07") to be written as genuine text instead of being auto-coerced to
# numbers by Excel's normal "smart" cell-entry parsing. Writing a formula
# whose result is a string, then copying *values only* into the real
# target cell, preserves the text without ever touching NumberFormat
# (which would otherwise register a brand-new style in the workbook).
$scratch = $totals.Range("Z1000")
function Set-TextValue($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# --- Step 1: relocate the existing "2022-Q1" sheet data to a new sheet ---
# Rename the current data sheet to "2022-Q3" (it will get fresh Q3 numbers
# below) and add a brand new sheet right after it that keeps the name
# "2022-Q1" and receives a copy of the original Q1 data + formatting.
$q1Sheet.Name = "2022-Q3"
$newQ1Sheet = $wb.Worksheets.Add($null, $q1Sheet)
$newQ1Sheet.Name = "2022-Q1"
$q1Sheet.Range("B1:H1").Copy($newQ1Sheet.Range("B1"))
$q1Sheet.Range("A2:H3").Copy($newQ1Sheet.Range("A2"))

# --- Step 2: overwrite the (renamed) sheet with the new 2022-Q3 data ---
$q3 = $q1Sheet

$q3.Range("A1:H3").ClearContents()

Set-TextValue $q3.Range("B1") "基金代码"
Set-TextValue $q3.Range("C1") "基金名称"
Set-TextValue $q3.Range("D1") "基金规模"
Set-TextValue $q3.Range("E1") "股票总仓位"
Set-TextValue $q3.Range("F1") "仓位占比"
Set-TextValue $q3.Range("G1") "持有市值(亿元)"
Set-TextValue $q3.Range("H1") "仓位排名"

$q3.Cells.Item(2, 1).Value = 0
Set-TextValue $q3.Range("B2") "001917"
Set-TextValue $q3.Range("C2") "招商量化精选股票A"
Set-TextValue $q3.Range("D2") "4.16"
Set-TextValue $q3.Range("E2") "92.70"
Set-TextValue $q3.Range("F2") "1.38"
Set-TextValue $q3.Range("G2") "0.0574"
$q3.Cells.Item(2, 8).Value = 9

$q3.Cells.Item(3, 1).Value = 1
Set-TextValue $q3.Range("B3") "007950"
Set-TextValue $q3.Range("C3") "招商量化精选股票C"
Set-TextValue $q3.Range("D3") "2.39"
Set-TextValue $q3.Range("E3") "92.70"
Set-TextValue $q3.Range("F3") "1.38"
Set-TextValue $q3.Range("G3") "0.0330"
$q3.Cells.Item(3, 8).Value = 9

# Re-apply the header / A-column style (s=2, same style used on the
# "总计" sheet's header row) now that ClearContents wiped formatting.
$totals.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$totals.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)
$q3.Range("A3").PasteSpecial(-4122)

# --- Step 3: update the "总计" sheet with the 2022-Q3 summary row + the
#     relocated 2022-Q1 summary row ---
$totals.Cells.Item(2, 2).Value = "2022-Q3"
$totals.Cells.Item(2, 3).Value = 2
$totals.Cells.Item(2, 4).Value = 0.09

$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$totals.Cells.Item(3, 1).Value = 1
$totals.Cells.Item(3, 2).Value = "2022-Q1"
$totals.Cells.Item(3, 3).Value = 2
$totals.Cells.Item(3, 4).Value = 0.6

# --- cleanup scratch area ---
$scratch.ClearContents()
